$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting all data rows (old 2..112) down by one
# (so old row 2 becomes row 3, ..., old row 112 becomes row 113).
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting from the header row; strip it back to
# the plain (unstyled) look used by every other data row.
$ws.Rows.Item(2).ClearFormats()

# New row 2 keeps the same constant columns as the rest of the table
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
# Calidad, Unidad de comercialización, Kg o Unidades, Clasificación)
# and gets fresh values for the variable columns.
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"
$ws.Range("D2").Value = 44860
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("N2").Value = "`$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 480
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
